# Kian Hou - Fixed the MIA 4
# Applies the data refresh for the new issuer (GREATOCEAN AUTOMOBILE SUPPLY
# SDN. BHD.) on the Knock-Out worksheet and clears the stale "S"/"U" helper
# columns (rows 11-35) that were left over from a previous calculation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header block: issuer name + CRA report date
# ---------------------------------------------------------------------
$ws.Range("E6").Value = "GREATOCEAN AUTOMOBILE SUPPLY SDN. BHD."
$ws.Range("J7").Value = "'2025-12-16"
$ws.Range("M7").Value = "GREATOCEAN AUTOMOBILE SUPPLY SDN. BHD."

# ---------------------------------------------------------------------
# Row 11 - Scoring by CRA Agency (Issuer's Credit Agency Score)
# ---------------------------------------------------------------------
$ws.Range("M11").Value = "'678"
$ws.Range("O11").Value = "'311"
$ws.Range("Q11").Value = "'533"
$ws.Range("S11").ClearContents()
$ws.Range("U11").ClearContents()

# ---------------------------------------------------------------------
# Row 12 - Scoring by CRA Agency (Credit Score Equivalent)
# ---------------------------------------------------------------------
$ws.Range("M12").Value = "B"
$ws.Range("O12").Value = "F"
$ws.Range("Q12").Value = "D"
$ws.Range("S12").ClearContents()
$ws.Range("U12").ClearContents()

# ---------------------------------------------------------------------
# Row 13 - Business operating since (years)
# ---------------------------------------------------------------------
$ws.Range("M13").Value = "'2011"

# ---------------------------------------------------------------------
# Row 18 - Credit Applications Approved for Last 12 months
# ---------------------------------------------------------------------
$ws.Range("M18").Value = "'4"
$ws.Range("O18").Value = "'0"
$ws.Range("Q18").Value = "'2"
$ws.Range("S18").ClearContents()
$ws.Range("U18").ClearContents()

# ---------------------------------------------------------------------
# Row 19 - Credit Applications Pending (no value change, only drop helpers)
# ---------------------------------------------------------------------
$ws.Range("S19").ClearContents()
$ws.Range("U19").ClearContents()

# ---------------------------------------------------------------------
# Row 20 - Legal Action taken (no value change, only drop helpers)
# ---------------------------------------------------------------------
$ws.Range("S20").ClearContents()
$ws.Range("U20").ClearContents()

# ---------------------------------------------------------------------
# Row 21 - Existing No. of Facility
# ---------------------------------------------------------------------
$ws.Range("M21").Value = "'16"
$ws.Range("O21").Value = "'10"
$ws.Range("Q21").Value = "'9"
$ws.Range("S21").ClearContents()
$ws.Range("U21").ClearContents()

# ---------------------------------------------------------------------
# Row 22 - Legal Suits
# ---------------------------------------------------------------------
$ws.Range("M22").Value = "'0"
$ws.Range("S22").ClearContents()
$ws.Range("U22").ClearContents()

# ---------------------------------------------------------------------
# Row 23 - Legal Case - Status
# ---------------------------------------------------------------------
$ws.Range("M23").Value = "No, No, No"
$ws.Range("S23").ClearContents()
$ws.Range("U23").ClearContents()

# ---------------------------------------------------------------------
# Row 24 - Trade / Credit Reference (cleared entirely, only drop helpers)
# ---------------------------------------------------------------------
$ws.Range("M24").ClearContents()
$ws.Range("O24").ClearContents()
$ws.Range("Q24").ClearContents()
$ws.Range("S24").ClearContents()
$ws.Range("U24").ClearContents()

# ---------------------------------------------------------------------
# Row 25 - Total Enquiries for Last 12 months
# ---------------------------------------------------------------------
$ws.Range("M25").Value = "'4"
$ws.Range("O25").Value = "'3"
$ws.Range("S25").ClearContents()
$ws.Range("U25").ClearContents()

# ---------------------------------------------------------------------
# Row 26 - Special Attention Account (no value change, only drop helpers)
# ---------------------------------------------------------------------
$ws.Range("S26").ClearContents()
$ws.Range("U26").ClearContents()

# ---------------------------------------------------------------------
# Row 27 - Summary of Total Liabilities (Outstanding)
# ---------------------------------------------------------------------
$ws.Range("M27").Value = "'15436493"
$ws.Range("O27").Value = "'1753657"
$ws.Range("Q27").Value = "'1723015"
$ws.Range("S27").ClearContents()
$ws.Range("U27").ClearContents()

# ---------------------------------------------------------------------
# Row 28 - Summary of Total Liabilities (Total Limit)
# ---------------------------------------------------------------------
$ws.Range("M28").Value = "'16364987"
$ws.Range("O28").Value = "'2146309"
$ws.Range("Q28").Value = "'2094126"
$ws.Range("S28").ClearContents()
$ws.Range("U28").ClearContents()

# ---------------------------------------------------------------------
# Row 29 - Overdraft facility outstanding (no value change, only drop helpers)
# ---------------------------------------------------------------------
$ws.Range("S29").ClearContents()
$ws.Range("U29").ClearContents()

# ---------------------------------------------------------------------
# Row 30 - Issuer's Total Banking Outstanding Facilities vs Limit
# ---------------------------------------------------------------------
$ws.Range("M30").Value = "YES, outstanding: 15520690.0, limit: 17714987.0"
$ws.Range("O30").Value = "YES, outstanding: 15520690.0, limit: 17714987.0"
$ws.Range("Q30").Value = "YES, outstanding: 15520690.0, limit: 17714987.0"
$ws.Range("S30").ClearContents()
$ws.Range("U30").ClearContents()

# ---------------------------------------------------------------------
# Row 31 - CCRIS Loan Account - Conduct Count
# ---------------------------------------------------------------------
$ws.Range("M31").Value = "current 1 month MIA1: 1, MIA2: 0, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 1, MIA2: 0, MIA3: 0, MIA4+: 0"
$ws.Range("O31").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 1 and /or past 6 months MIA1: 1, MIA2: 1, MIA3: 1, MIA4+: 2"
$ws.Range("Q31").Value = "current 1 month MIA1: 0, MIA2: 1, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 6, MIA2: 3, MIA3: 1, MIA4+: 0"
$ws.Range("S31").ClearContents()
$ws.Range("U31").ClearContents()

# ---------------------------------------------------------------------
# Row 33 - Issuer's Total Non-Bank Lender Outstanding vs Limit (helpers only)
# ---------------------------------------------------------------------
$ws.Range("S33").ClearContents()
$ws.Range("U33").ClearContents()

# ---------------------------------------------------------------------
# Row 34 - Non-Bank Lender Credit Information (NLCI) - Conduct Count
# ---------------------------------------------------------------------
$ws.Range("M34").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0"
$ws.Range("O34").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0"
$ws.Range("Q34").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0"
$ws.Range("S34").ClearContents()
$ws.Range("U34").ClearContents()

# ---------------------------------------------------------------------
# Row 35 - Non-Bank Lender Credit Information (NLCI) - Legal Status
# ---------------------------------------------------------------------
$ws.Range("M35").Value = "WITHDRAWN"
$ws.Range("O35").Value = "WITHDRAWN"
$ws.Range("Q35").Value = "WITHDRAWN"
$ws.Range("S35").ClearContents()
$ws.Range("U35").ClearContents()
